$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply each cell update as text, preserving the original inline-string / text
# cell type (these look numeric so Excel would otherwise auto-convert them to
# real numbers). We force a text number format, assign the value, then restore
# the default "Normal" style so no stray formatting is left behind.
function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "26.273.92"
Set-TextValue "E2" "  +0.63%  "
Set-TextValue "D3" "1.665.79"
Set-TextValue "D5" "218.75"
Set-TextValue "E5" "  +0.57%  "
Set-TextValue "D6" "0.5324"
Set-TextValue "E6" "  +1.52%  "
Set-TextValue "E7" "  +0.75%  "
Set-TextValue "E8" "  +1.27%  "
Set-TextValue "E9" "  +0.55%  "
Set-TextValue "E10" "  +0.95%  "
Set-TextValue "D11" "0.07819"
Set-TextValue "E11" "  +0.20%  "
Set-TextValue "D12" "4.558"
Set-TextValue "E12" "  +1.28%  "
Set-TextValue "D13" "1.670.49"
Set-TextValue "E13" "  +1.11%  "
Set-TextValue "D14" "1.894.36"
Set-TextValue "D15" "0.5527"
Set-TextValue "E15" "  +1.07%  "
Set-TextValue "D16" "0.0₅8218"
Set-TextValue "E16" "  +0.30%  "
Set-TextValue "D17" "65.75"
Set-TextValue "E18" "  +0.79%  "
Set-TextValue "D19" "4.675"
Set-TextValue "E19" "  +2.10%  "
Set-TextValue "D20" "193.61"
Set-TextValue "E20" "  +1.25%  "
Set-TextValue "E21" "  +1.56%  "
Set-TextValue "D22" "6.034"
Set-TextValue "E23" "  +0.74%  "
Set-TextValue "D24" "145.58"
Set-TextValue "E24" "  +2.54%  "
Set-TextValue "E25" "  -0.82%  "
Set-TextValue "D26" "7.199"
Set-TextValue "E26" "  -0.50%  "
Set-TextValue "D27" "16.17"
Set-TextValue "E27" "  +0.19%  "
Set-TextValue "D28" "1.485"
Set-TextValue "E28" "  +3.70%  "
Set-TextValue "D29" "0.05897"
Set-TextValue "E29" "  +0.01%  "
Set-TextValue "D30" "1.282"
Set-TextValue "E30" "  +0.13%  "
Set-TextValue "D31" "3.597"
Set-TextValue "E31" "  +2.29%  "
Set-TextValue "E32" "  +0.88%  "
Set-TextValue "E33" "  +1.63%  "
Set-TextValue "D34" "0.9627"
Set-TextValue "E34" "  +1.28%  "
Set-TextValue "E35" "  +1.41%  "
Set-TextValue "E36" "  +0.34%  "
Set-TextValue "D37" "0.5795"
Set-TextValue "E37" "  +2.18%  "
Set-TextValue "D38" "0.01609"
Set-TextValue "E38" "  -0.47%  "
Set-TextValue "D39" "0.8645"
Set-TextValue "E39" "  +1.93%  "
Set-TextValue "D40" "5.836"
Set-TextValue "E40" "  +0.24%  "
Set-TextValue "B41" "Maker"
Set-TextValue "C41" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D41" "1.048.71"
Set-TextValue "E41" "  +1.91%  "
Set-TextValue "B42" "PaxDollar"
Set-TextValue "C42" "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue "D42" "1.009"
Set-TextValue "E42" "  +0.69%  "
Set-TextValue "D43" "104.01"
Set-TextValue "E43" "  +1.31%  "
Set-TextValue "D44" "1.805.05"
Set-TextValue "E44" "  +0.40%  "
Set-TextValue "D45" "57.68"
Set-TextValue "D46" "1.013"
Set-TextValue "E46" "  +0.96%  "
Set-TextValue "E47" "  -5.69%  "
Set-TextValue "B48" "EnergySwap"
Set-TextValue "C48" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D48" "8.080"
Set-TextValue "E48" "  +2.72%  "
Set-TextValue "B49" "Mantle"
Set-TextValue "C49" "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D49" "0.4379"
Set-TextValue "E49" "  +1.71%  "
Set-TextValue "D50" "0.05162"
Set-TextValue "E50" "  -0.09%  "
Set-TextValue "E51" "  -3.07%  "
